# Apply the "Added vanilla needy modules" edit to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "modules-config-details"

# Append the three new needy-module rows (110-112), following the same
# column layout used by the existing rows: Name, CodeName, Category,
# PDF path, Author, Release date, Flag.
$newRows = @(
    @("Capacitor Discharge", "NeedyCapacitor", 1, "modules/Capacitor Discharge.pdf", "Steel Crate Games", "2015-10-08", 2),
    @("Knob",                "NeedyKnob",      1, "modules/Knob.pdf",                "Steel Crate Games", "2015-10-08", 2),
    @("Venting Gas",         "NeedyVentGas",   1, "modules/Venting Gas.pdf",          "Steel Crate Games", "2015-10-08", 2)
)

$row = 110
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $row++
}

# Update the view/selection to match where the author ended up.
$null = $ws.Range("J99").Select()
